# Generate Report for Archive
#
# 1) The localization status changes from "Ready for handoff" to
#    "In Translation" on every sheet that shows it (Overview!E2:F2,
#    zh-cn!C2, de-de!C2).
# 2) The two "date/status" columns that used to be sized for the longer
#    "Ready for handoff" text (stored width 17.2159881591797) are
#    narrowed to fit "In Translation" (stored width 13.4101845877511):
#    Overview columns E & F, and column C on both the zh-cn and de-de
#    sheets.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# New narrower width. (ColumnWidth is quantized to whole pixels by the
# engine, so 12.5 is the value that lands closest to the target stored
# width of 13.4101845877511.)
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
if ($wsOverview.Range("E2").Value2 -eq $oldStatus) { $wsOverview.Range("E2").Value = $newStatus }
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) { $wsOverview.Range("F2").Value = $newStatus }
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
if ($wsZhCn.Range("C2").Value2 -eq $oldStatus) { $wsZhCn.Range("C2").Value = $newStatus }
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
if ($wsDeDe.Range("C2").Value2 -eq $oldStatus) { $wsDeDe.Range("C2").Value = $newStatus }
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
